$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date and status values in row 4
$ws.Range("A4").Value = "29/08/2016"
$ws.Range("H4").Value = "Done"
$ws.Range("I4").Value = "Done"
$ws.Range("J4").Value = "In Progress"
$ws.Range("K4").Value = "In Progress"
$ws.Range("L4").Value = "Done"

# Update the view: scroll so column I is the top-left visible column,
# and select L4 as the active cell
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("L4").Select()
